$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 511, shifting existing rows 511-548 down to 513-550.
$ws.Rows("511:512").Insert()

# Row 511 - new record (Piña, Primera)
$ws.Range("A511").Value = 10
$ws.Range("B511").Value = "Vega Modelo de Temuco"
$ws.Range("C511").Value = "La Araucanía"
$ws.Range("D511").Value = 44826
$ws.Range("E511").Value = 9
$ws.Range("F511").Value = "Fruta"
$ws.Range("G511").Value = 100108
$ws.Range("H511").Value = "Tropicales y subtropicales"
$ws.Range("I511").Value = 100108005
$ws.Range("J511").Value = "Piña"
$ws.Range("K511").Value = "Caramelo"
$ws.Range("L511").Value = "Primera"
$ws.Range("M511").Value = 90
$ws.Range("N511").Value = 22000
$ws.Range("O511").Value = 23000
$ws.Range("P511").Value = 22611
$ws.Range("Q511").Value = "$/caja 12 unidades"
$ws.Range("R511").Value = "Ecuador"
$ws.Range("S511").Value = 1884
$ws.Range("T511").Value = 12

# Row 512 - new record (Piña, Segunda)
$ws.Range("A512").Value = 10
$ws.Range("B512").Value = "Vega Modelo de Temuco"
$ws.Range("C512").Value = "La Araucanía"
$ws.Range("D512").Value = 44826
$ws.Range("E512").Value = 9
$ws.Range("F512").Value = "Fruta"
$ws.Range("G512").Value = 100108
$ws.Range("H512").Value = "Tropicales y subtropicales"
$ws.Range("I512").Value = 100108005
$ws.Range("J512").Value = "Piña"
$ws.Range("K512").Value = "Caramelo"
$ws.Range("L512").Value = "Segunda"
$ws.Range("M512").Value = 125
$ws.Range("N512").Value = 22000
$ws.Range("O512").Value = 23000
$ws.Range("P512").Value = 22480
$ws.Range("Q512").Value = "$/caja 14 unidades"
$ws.Range("R512").Value = "Ecuador"
$ws.Range("S512").Value = 1606
$ws.Range("T512").Value = 14
